$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row (2-203).
# The value changes from 45179 (2023-09-10) to 45180 (2023-09-11) for all rows.
$ws.Range("C2:C203").Value = 45180
